# Update "想去人数" (F column) counts across the sheets of the workbook.
# Sheet order (per workbook.xml): 1=展览, 2=演出, 3=本地生活, 4=全部类型

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 413
$ws1.Range("F4").Value = 1170
$ws1.Range("F5").Value = 49
$ws1.Range("F7").Value = 37
$ws1.Range("F9").Value = 524
$ws1.Range("F10").Value = 373
$ws1.Range("F11").Value = 429
$ws1.Range("F14").Value = 365
$ws1.Range("F15").Value = 42
$ws1.Range("F16").Value = 66
$ws1.Range("F18").Value = 550
$ws1.Range("F19").Value = 1466
$ws1.Range("F20").Value = 5707
$ws1.Range("F21").Value = 92
$ws1.Range("F22").Value = 1601
$ws1.Range("F23").Value = 381
$ws1.Range("F24").Value = 53
$ws1.Range("F25").Value = 29
$ws1.Range("F26").Value = 5287
$ws1.Range("F27").Value = 5287
$ws1.Range("F28").Value = 130
$ws1.Range("F29").Value = 86
$ws1.Range("F30").Value = 1538
$ws1.Range("F32").Value = 26
$ws1.Range("F33").Value = 50
$ws1.Range("F34").Value = 36
$ws1.Range("F35").Value = 667
$ws1.Range("F36").Value = 100

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value = 28
$ws2.Range("F5").Value = 158
$ws2.Range("F8").Value = 164
$ws2.Range("F11").Value = 4

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 9410

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 9410
$ws4.Range("F6").Value = 413
$ws4.Range("F7").Value = 1170
$ws4.Range("F8").Value = 49
$ws4.Range("F10").Value = 37
$ws4.Range("F12").Value = 373
$ws4.Range("F13").Value = 429
$ws4.Range("F15").Value = 365
$ws4.Range("F16").Value = 42
$ws4.Range("F17").Value = 66
$ws4.Range("F21").Value = 550
$ws4.Range("F22").Value = 1466
$ws4.Range("F23").Value = 5707
$ws4.Range("F24").Value = 92
$ws4.Range("F25").Value = 1601
$ws4.Range("F28").Value = 381
$ws4.Range("F29").Value = 4
$ws4.Range("F31").Value = 5287
$ws4.Range("F32").Value = 5287
$ws4.Range("F33").Value = 130
$ws4.Range("F34").Value = 86
$ws4.Range("F35").Value = 1538
$ws4.Range("F37").Value = 26
$ws4.Range("F38").Value = 667
$ws4.Range("F39").Value = 100
